$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "69.505.53"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +0.25%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.676.09"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("E4").Value = "  +0.01%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "654.37"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -3.50%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "159.58"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("E7").Value = "  +0.04%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.496"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +0.46%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.144"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -1.12%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "7.05"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -0.28%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.439"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +0.70%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.0000231"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -0.36%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "4.296.14"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -0.28%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "32.54"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +0.56%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "3.683.13"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -0.32%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "69.505.74"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("E17").Value = "  +1.47%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "15.98"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -0.33%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "6.42"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +0.00%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "466.56"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.40%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "9.71"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -2.59%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.640"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -1.74%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "79.59"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.35%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "3.824.18"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -0.21%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "0.0000125"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +1.04%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "10.86"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.71%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "8.93"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -1.99%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.63"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -1.87%  "
$ws.Range("E30").Value = "  -5.48%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("E32").Value = "  +0.11%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "26.63"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -1.22%  "
$ws.Range("E34").Value = "  -2.64%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.163"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +3.07%  "
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "3.667.16"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -0.20%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "8.38"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +1.59%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "5.99"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -4.16%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "178.93"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +4.97%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -0.06%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "2.20"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -1.60%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.0893"
$cell.Style = "Normal"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.929"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -1.41%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "46.82"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -1.80%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "2.72"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +0.14%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "1.28"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +0.04%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "27.21"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -3.49%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.000267"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -4.46%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "7.79"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("E51").Value = "  -3.78%  "
